{"js": "// Add a new paragraph \"Add new line\" at the end of the document body,\n// right after the existing \"For test\" paragraph (matches the diff: a\n// new <w:p> with an en-US run is appended before the sectPr).\nconst body = context.document.body;\nbody.insertParagraph(\"Add new line\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Add a new paragraph \"Add new line\" at the end of the document,\n# right after the existing \"For test\" paragraph (matches the diff: a\n# new <w:p> with an en-US run is appended before the sectPr).\n$d = $word.ActiveDocument\n$p = $d.Paragraphs.Add()\n$p.Range.Text = \"Add new line\"\n"}
